# Automatic update of files.
# Update the "Förändrad" (C) date-serial for every data row, and re-shuffle
# the data (Beteckning/Datum/Area) among rows 4,5,6,7,9,10 as per source
# update, leaving row 8 (and rows 2,3) values in A/B/G untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumps from 46064 to 46065 for every data row (2-10)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 46065
}

# New values for columns A (Beteckning), B (Datum), G (Area (ha)) for rows
# that get reshuffled.
$ws.Cells.Item(4, 1).Value = "A 14516-2023"
$ws.Cells.Item(4, 2).Value = 45012.86600694444
$ws.Cells.Item(4, 7).Value = 0.4

$ws.Cells.Item(5, 1).Value = "A 26262-2024"
$ws.Cells.Item(5, 2).Value = 45468.66077546297
$ws.Cells.Item(5, 7).Value = 0.6

$ws.Cells.Item(6, 1).Value = "A 4156-2023"
$ws.Cells.Item(6, 2).Value = 44953
$ws.Cells.Item(6, 7).Value = 1.5

$ws.Cells.Item(7, 1).Value = "A 14517-2023"
$ws.Cells.Item(7, 2).Value = 45012
$ws.Cells.Item(7, 7).Value = 0.6

$ws.Cells.Item(9, 1).Value = "A 4159-2023"
$ws.Cells.Item(9, 2).Value = 44953
$ws.Cells.Item(9, 7).Value = 0.5

$ws.Cells.Item(10, 1).Value = "A 50762-2025"
$ws.Cells.Item(10, 2).Value = 45946
$ws.Cells.Item(10, 7).Value = 2.7
